$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add new header values in P1 and Q1
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: swap I<->K and M<->O values, and add new P/Q columns with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I column: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K column: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M column: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O column: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P column: new, value 2
    $ws.Cells.Item($r, 17).Value = 2  # Q column: new, value 2
}
